# Avoid recalculating the existing SUM/AVERAGE formulas in row 10; only the
# raw data cells themselves are being updated (matches source diff, which
# leaves the cached formula results in row 10 untouched).
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the accuracy summary text
$ws.Range("A11").Value = "Akurasi Pengujian = 32.73%"

# Row 3 (Nomor Telepon indicator)
$ws.Range("H3").Value = 1001
$ws.Range("I3").Value = 0.8833786231884058

# Row 4 (Riwayat Pendidikan indicator)
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 1

# Row 5 (Alamat Institusi / Catatan Rekam Medis indicators)
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("E5").Value = 433
$ws.Range("F5").Value = 1

# Row 8 (No HP Siswa indicator)
$ws.Range("H8").Value = 472
$ws.Range("I8").Value = 0.90625
